# GNTX Yearly Financials - "Doing Updates for Financials"
#
# The source data adds a new fiscal-year column (period ending 2018-12-31,
# serial 43465) as the new column D on the Income Statement, Balance Sheet
# and Cash Flow Statement blocks. All previously existing year columns
# (old D:K) shift one column to the right (new E:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlShiftToRight = -4161 ; xlPasteFormats = -4122
$xlShiftToRight = -4161
$xlPasteFormats = -4122

# Insert a blank column in front of the existing "D" column, for every row
# that actually holds the year-over-year tables (rows 5-102). Using a
# bounded range (not whole-column) keeps blank separator rows (5,6,37,79)
# untouched and avoids materialising the fully-empty rows 36/78.
$ws.Range("D5:D102").Insert($xlShiftToRight)

# Copy the number formats/fonts/alignment from the (now shifted) old "D"
# column - which landed in "E" - back onto the freshly inserted "D" column,
# one contiguous data block at a time so the blank separator rows keep no
# cell at all (matching the original workbook's sparse rows).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial($xlPasteFormats)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial($xlPasteFormats)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---- Income Statement (new FY2018 column) ----
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1834100
$ws.Range("D9").Value = 1143600
$ws.Range("D10").Value = 690500
$ws.Range("D12").Value = 107100
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1325900
$ws.Range("D18").Value = 508100
$ws.Range("D20").Value = 14700
$ws.Range("D21").Value = 625000
$ws.Range("D22").Value = 800
$ws.Range("D23").Value = 522000
$ws.Range("D24").Value = 84200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 437900
$ws.Range("D27").Value = 437900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -14700
$ws.Range("D33").Value = 437900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 437900

# ---- Balance Sheet (new FY2018 column) ----
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 217000
$ws.Range("D42").Value = 169400
$ws.Range("D43").Value = 213500
$ws.Range("D44").Value = 225300
$ws.Range("D45").Value = 25700
$ws.Range("D46").Value = 850900
$ws.Range("D47").Value = 138000
$ws.Range("D48").Value = 498500
$ws.Range("D49").Value = 598100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2085400
$ws.Range("D57").Value = 92800
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 76400
$ws.Range("D60").Value = 169200
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 54500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 223700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1102500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1861800
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement (new FY2018 column) ----
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 437900
$ws.Range("D83").Value = 102200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 552400
$ws.Range("D91").Value = -91600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -185800
$ws.Range("D96").Value = -116600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -719300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -352700
